# Insert two new weekly-report rows for "Pepino ensalada" (Terminal La Palmera
# de La Serena), shifting the existing rows 245:340 down to 247:342.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 245 - this pushes every row that
# was 245..340 down to 247..342, inheriting formatting (incl. the date-style
# column D) from the row above, same as a real Excel "Insert Copied Cells"
# done from the UI.
$ws.Rows("245:246").Insert()

# ---- New row 245 --------------------------------------------------------
$ws.Cells.Item(245, 1).Value = 8
$ws.Cells.Item(245, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(245, 3).Value = "Coquimbo"
$ws.Cells.Item(245, 4).Value = 44468
$ws.Cells.Item(245, 5).Value = 4
$ws.Cells.Item(245, 6).Value = 100112043
$ws.Cells.Item(245, 7).Value = "Pepino ensalada"
$ws.Cells.Item(245, 8).Value = "Sin especificar"
$ws.Cells.Item(245, 9).Value = "Primera"
$ws.Cells.Item(245, 10).Value = 800
$ws.Cells.Item(245, 11).Value = 15000
$ws.Cells.Item(245, 12).Value = 16000
$ws.Cells.Item(245, 13).Value = 15500
$ws.Cells.Item(245, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(245, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(245, 16).Value = 258
$ws.Cells.Item(245, 17).Value = 60
$ws.Cells.Item(245, 18).Value = "Hortaliza"

# ---- New row 246 --------------------------------------------------------
$ws.Cells.Item(246, 1).Value = 8
$ws.Cells.Item(246, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(246, 3).Value = "Coquimbo"
$ws.Cells.Item(246, 4).Value = 44468
$ws.Cells.Item(246, 5).Value = 4
$ws.Cells.Item(246, 6).Value = 100112043
$ws.Cells.Item(246, 7).Value = "Pepino ensalada"
$ws.Cells.Item(246, 8).Value = "Sin especificar"
$ws.Cells.Item(246, 9).Value = "Segunda"
$ws.Cells.Item(246, 10).Value = 560
$ws.Cells.Item(246, 11).Value = 11000
$ws.Cells.Item(246, 12).Value = 12000
$ws.Cells.Item(246, 13).Value = 11500
$ws.Cells.Item(246, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(246, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(246, 16).Value = 115
$ws.Cells.Item(246, 17).Value = 100
$ws.Cells.Item(246, 18).Value = "Hortaliza"
